# Refresh the crypto price/volume snapshot values (mirrors the scheduled
# GitHub Actions scraper commit). Coin (B) and Link (C) cells are plain
# text already; Price (D) cells are numeric-looking text, so we pin their
# NumberFormat to "@" (Text) before writing so Excel keeps them as strings
# (otherwise e.g. "1.00" would silently become the number 1). Volume (E)
# cells already contain "%" / spaces, so plain string assignment is safe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.430.05'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.275.68'
$ws.Range("E3").Value = '  +3.33%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.86'
$ws.Range("E5").Value = '  +2.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.56'
$ws.Range("E6").Value = '  +3.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.274.82'
$ws.Range("E8").Value = '  +3.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  +3.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.78'
$ws.Range("E11").Value = '  +2.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.497'
$ws.Range("E12").Value = '  -3.53%  '
$ws.Range("E13").Value = '  +3.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.09'
$ws.Range("E14").Value = '  +2.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.817.82'
$ws.Range("E15").Value = '  +3.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.541.41'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.42'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.281.61'
$ws.Range("E18").Value = '  +3.54%  '
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '505.70'
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.57'
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.756'
$ws.Range("E22").Value = '  +4.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.17'
$ws.Range("E23").Value = '  +2.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.69'
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.84'
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +2.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.28'
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.42'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.128'
$ws.Range("E30").Value = '  +46.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.01'
$ws.Range("E31").Value = '  -2.21%  '
$ws.Range("E32").Value = '  -3.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.14'
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.50'
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.51'
$ws.Range("E37").Value = '  +23.49%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0793'
$ws.Range("E38").Value = '  +17.21%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '55.54'
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '499.76'
$ws.Range("E40").Value = '  -1.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0425'
$ws.Range("E41").Value = '  +2.14%  '
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.82'
$ws.Range("E43").Value = '  +1.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.62'
$ws.Range("E44").Value = '  +7.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.295'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.012.86'
$ws.Range("E46").Value = '  +6.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '29.20'
$ws.Range("E47").Value = '  +5.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.48'
$ws.Range("E48").Value = '  +5.58%  '
$ws.Range("E49").Value = '  +2.80%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.39'
$ws.Range("E51").Value = '  -0.25%  '
